# Fruta / hortaliza, semanal
# Re-order the weekly price rows (2-12) for columns D (Fecha), J (Volumen),
# K (Precio mínimo), L (Precio máximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) according to the new weekly sequence, while leaving all
# other columns (which are identical across rows) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for the columns that move, keyed by row number,
# before any writes happen (so later writes don't clobber values we still
# need to read).
$cols = @("D", "J", "K", "L", "M", "P")
$original = @{}
for ($r = 2; $r -le 12; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $original[$r] = $rowVals
}

# Mapping of destination row -> source row (where the data should come from)
$mapping = @{
    2  = 7
    3  = 5
    4  = 8
    5  = 3
    6  = 9
    7  = 12
    8  = 4
    9  = 2
    10 = 11
    11 = 6
    12 = 10
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
